# Auto-generated Excel COM-interop script to apply market price / profit updates
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8825
$ws.Range("I51").Value = 24999
$ws.Range("K51").Value = 24999
$ws.Range("M51").Value = -24515
$ws.Range("H62").Value = 8463.6
$ws.Range("I62").Value = 7517.625
$ws.Range("K62").Value = 7517.625
$ws.Range("M62").Value = -6893.625
$ws.Range("H65").Value = 8463.6
$ws.Range("I65").Value = 7517.625
$ws.Range("K65").Value = 37588.125
$ws.Range("M65").Value = -34468.125
$ws.Range("H98").Value = 1458.2941
$ws.Range("I98").Value = 1233.2188
$ws.Range("K98").Value = 1233.2188
$ws.Range("M98").Value = 264.7811999999999
$ws.Range("H122").Value = 1458.2941
$ws.Range("I122").Value = 1233.2188
$ws.Range("K122").Value = 3699.6564
$ws.Range("M122").Value = -1249.6564
$ws.Range("H132").Value = 10529.75
$ws.Range("I132").Value = 10529.75
$ws.Range("K132").Value = 31589.25
$ws.Range("M132").Value = -29059.25
$ws.Range("H138").Value = 3924.1
$ws.Range("J138").Value = 6608.1577
$ws.Range("L138").Value = 19824.4731
$ws.Range("N138").Value = -30104.4731
$ws.Range("H140").Value = 402164.66
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8459.097
$ws.Range("I32").Value = 6718.2915
$ws.Range("J32").Value = 14427.571
$ws.Range("K32").Value = 6718.2915
$ws.Range("L32").Value = 14427.571
$ws.Range("M32").Value = -6431.2915
$ws.Range("N32").Value = -15001.571
$ws.Range("H61").Value = 3952.682
$ws.Range("I61").Value = 3831.3333
$ws.Range("K61").Value = 3831.3333
$ws.Range("M61").Value = -3619.3333
$ws.Range("H74").Value = 4183.273
$ws.Range("I74").Value = 4101.6
$ws.Range("K74").Value = 4101.6
$ws.Range("M74").Value = -3227.6
$ws.Range("H77").Value = 4183.273
$ws.Range("I77").Value = 4101.6
$ws.Range("K77").Value = 20508
$ws.Range("M77").Value = -16140
$ws.Range("H97").Value = 1438.0588
$ws.Range("I97").Value = 1142
$ws.Range("K97").Value = 1142
$ws.Range("M97").Value = -646
$ws.Range("H110").Value = 2981.5
$ws.Range("I110").Value = 2725.3333
$ws.Range("K110").Value = 2725.3333
$ws.Range("M110").Value = -680.3332999999998
$ws.Range("H132").Value = 2663.762
$ws.Range("I132").Value = 2318.6765
$ws.Range("K132").Value = 6956.029500000001
$ws.Range("M132").Value = -4426.029500000001
$ws.Range("H136").Value = 3952.682
$ws.Range("I136").Value = 3831.3333
$ws.Range("K136").Value = 11493.9999
$ws.Range("M136").Value = -8943.999899999999

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2675.182
$ws.Range("I105").Value = 2319.25
$ws.Range("J105").Value = 3624.3333
$ws.Range("K105").Value = 2319.25
$ws.Range("L105").Value = 3624.3333
$ws.Range("M105").Value = -572.25
$ws.Range("N105").Value = -7118.3333
$ws.Range("H107").Value = 1316.5714
$ws.Range("I107").Value = 1165.8158
$ws.Range("K107").Value = 1165.8158
$ws.Range("M107").Value = 754.1841999999999

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 755000.5
$ws.Range("J4").Value = 1010001
$ws.Range("L4").Value = 1010001
$ws.Range("N4").Value = -1010225
$ws.Range("H132").Value = 2558.6667
$ws.Range("J132").Value = 1474.5
$ws.Range("L132").Value = 4423.5
$ws.Range("N132").Value = -9483.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 700.5
$ws.Range("I25").Value = 1001
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 3003
$ws.Range("L25").Value = 1200
$ws.Range("M25").Value = -2834
$ws.Range("N25").Value = -1538
$ws.Range("H30").Value = 700.5
$ws.Range("I30").Value = 1001
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 3003
$ws.Range("L30").Value = 1200
$ws.Range("M30").Value = -2901
$ws.Range("N30").Value = -1404
$ws.Range("H39").Value = 4987.7144
$ws.Range("J39").Value = 4987.7144
$ws.Range("L39").Value = 14963.1432
$ws.Range("N39").Value = -15551.1432
$ws.Range("H109").Value = 1968
$ws.Range("J109").Value = 2974
$ws.Range("L109").Value = 8922
$ws.Range("N109").Value = -11002
$ws.Range("H119").Value = 2860.5715
$ws.Range("I119").Value = 2670.6667
$ws.Range("J119").Value = 4000
$ws.Range("K119").Value = 8012.000100000001
$ws.Range("L119").Value = 12000
$ws.Range("M119").Value = -3174.000100000001
$ws.Range("N119").Value = -21676

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9982.833000000001
$ws.Range("I5").Value = 9982.833000000001
$ws.Range("K5").Value = 9982.833000000001
$ws.Range("M5").Value = -9870.833000000001
$ws.Range("H80").Value = 2692.1
$ws.Range("I80").Value = 2264.6667
$ws.Range("J80").Value = 3333.25
$ws.Range("K80").Value = 2264.6667
$ws.Range("L80").Value = 3333.25
$ws.Range("M80").Value = -1266.6667
$ws.Range("N80").Value = -5329.25
$ws.Range("H83").Value = 2692.1
$ws.Range("I83").Value = 2264.6667
$ws.Range("J83").Value = 3333.25
$ws.Range("K83").Value = 11323.3335
$ws.Range("L83").Value = 16666.25
$ws.Range("M83").Value = -6331.333500000001
$ws.Range("N83").Value = -26650.25
$ws.Range("H113").Value = 1202.7142
$ws.Range("I113").Value = 1499
$ws.Range("K113").Value = 1499
$ws.Range("M113").Value = 671

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 27006.334
$ws.Range("I5").Value = 12009
$ws.Range("K5").Value = 12009
$ws.Range("M5").Value = -11896
$ws.Range("H7").Value = 9286.909
$ws.Range("I7").Value = 9959.823
$ws.Range("K7").Value = 9959.823
$ws.Range("M7").Value = -9847.823
$ws.Range("H125").Value = 124994.5
$ws.Range("J125").Value = 124994.5
$ws.Range("L125").Value = 124994.5
$ws.Range("N125").Value = -134834.5
$ws.Range("H126").Value = 9286.909
$ws.Range("I126").Value = 9959.823
$ws.Range("K126").Value = 29879.469
$ws.Range("M126").Value = -27409.469
$ws.Range("H132").Value = 3062.5715
$ws.Range("I132").Value = 1996
$ws.Range("K132").Value = 5988
$ws.Range("M132").Value = -3458
$ws.Range("H136").Value = 4260.8076
$ws.Range("I136").Value = 3265.6667
$ws.Range("K136").Value = 9797.000100000001
$ws.Range("M136").Value = -7247.000100000001

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 24512.25
$ws.Range("I21").Value = 19007.5
$ws.Range("K21").Value = 19007.5
$ws.Range("M21").Value = -18772.5
$ws.Range("H30").Value = 30009
$ws.Range("I30").Value = 30009
$ws.Range("K30").Value = 30009
$ws.Range("M30").Value = -29902
$ws.Range("H35").Value = 24512.25
$ws.Range("I35").Value = 19007.5
$ws.Range("K35").Value = 19007.5
$ws.Range("M35").Value = -18717.5
$ws.Range("H64").Value = 99967.25
$ws.Range("J64").Value = 99967.25
$ws.Range("L64").Value = 99967.25
$ws.Range("N64").Value = -100463.25
$ws.Range("H67").Value = 99967.25
$ws.Range("J67").Value = 99967.25
$ws.Range("L67").Value = 99967.25
$ws.Range("N67").Value = -101683.25
$ws.Range("H81").Value = 2880.4167
$ws.Range("I81").Value = 1413
$ws.Range("J81").Value = 3928.5715
$ws.Range("K81").Value = 2826
$ws.Range("L81").Value = 7857.143
$ws.Range("M81").Value = -1765
$ws.Range("N81").Value = -9979.143
$ws.Range("H84").Value = 2880.4167
$ws.Range("I84").Value = 1413
$ws.Range("J84").Value = 3928.5715
$ws.Range("K84").Value = 14130
$ws.Range("L84").Value = 39285.715
$ws.Range("M84").Value = -8826
$ws.Range("N84").Value = -49893.715
$ws.Range("H122").Value = 7396.0835
$ws.Range("J122").Value = 11458.833
$ws.Range("L122").Value = 34376.499
$ws.Range("N122").Value = -39276.499
$ws.Range("H126").Value = 6829
$ws.Range("I126").Value = 7124.5625
$ws.Range("K126").Value = 21373.6875
$ws.Range("M126").Value = -18903.6875
$ws.Range("H128").Value = 16784764
$ws.Range("J128").Value = 16784764
$ws.Range("L128").Value = 16784764
$ws.Range("N128").Value = -16794724
